$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.503.84"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.62"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.07"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4692"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3978"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.27"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -11.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08056"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.021"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.82"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.958"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.194"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.92"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06560"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.33"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.514"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.513.61"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.306"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.085.77"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.04"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.28"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.086"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.569"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.42"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09483"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9563"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.472"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.603"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.304"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06099"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02251"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.220"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.152"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5995"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1894"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.31"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.266"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5676"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.25"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.411"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.941"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06788"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.23"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
